$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# Add the two new keyword rows (order matters: it controls the order in
# which new entries are appended to the shared-strings table).
$ws.Range("A3").Value = "Some spaces"
$ws.Range("B3").Value = "Default"
$ws.Range("B4").Value = "Income"
$ws.Range("A4").Value = "`$tr@ng€ Нейм"

# Move the active selection down to the new last row, like a user would
# after typing the new data.
[void]$ws.Range("A4").Select()

# Match the page setup (paper size / orientation) recorded for the sheet.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
